$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 400m (M) results: Winner, 2nd Place, 3rd Place
$ws.Range("B6").Value = "Quincy Hall"
$ws.Range("C6").Value = "Mattew Hudson-Smith"
$ws.Range("D6").Value = "Muzala Samukonga"

# Normalize formatting across the whole table: remove the bold/10pt header
# font, bring everything to the default 11pt font, and center all columns
# horizontally.
$tbl = $ws.Range("A1:D9")
$tbl.Font.Bold = $false
$tbl.Font.Size = 11
$tbl.HorizontalAlignment = -4108

# Let the (now wider) columns re-fit their contents.
$tbl.EntireColumn.AutoFit()

$ws.Range("D7").Select()
